$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '35.281.85'
$c.Style = "Normal"
$ws.Range('E2').Value = '  -0.48%  '
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '1.913.01'
$c.Style = "Normal"
$ws.Range('E3').Value = '  -0.09%  '
$ws.Range('E4').Value = '  +0.42%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '0.724'
$c.Style = "Normal"
$ws.Range('E5').Value = '  +9.12%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '255.81'
$c.Style = "Normal"
$ws.Range('E6').Value = '  +3.61%  '
$ws.Range('E7').Value = '  +0.38%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '41.94'
$c.Style = "Normal"
$ws.Range('E8').Value = '  +0.34%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.368'
$c.Style = "Normal"
$ws.Range('E9').Value = '  +5.76%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '53.39'
$c.Style = "Normal"
$ws.Range('E10').Value = '  +0.32%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.0764'
$c.Style = "Normal"
$ws.Range('E11').Value = '  +6.05%  '
$ws.Range('E12').Value = '  -0.15%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '13.13'
$c.Style = "Normal"
$ws.Range('E13').Value = '  +6.80%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '2.195.69'
$c.Style = "Normal"
$ws.Range('E14').Value = '  +0.49%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '0.737'
$c.Style = "Normal"
$ws.Range('E15').Value = '  +5.24%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '5.00'
$c.Style = "Normal"
$ws.Range('E16').Value = '  +3.40%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '1.936.35'
$c.Style = "Normal"
$ws.Range('E17').Value = '  +1.02%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '35.312.92'
$c.Style = "Normal"
$ws.Range('E18').Value = '  -0.35%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '75.26'
$c.Style = "Normal"
$ws.Range('E19').Value = '  +4.17%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '0.0₃0849'
$c.Style = "Normal"
$ws.Range('E20').Value = '  +3.43%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '246.20'
$c.Style = "Normal"
$ws.Range('E21').Value = '  +1.92%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '13.11'
$c.Style = "Normal"
$ws.Range('E22').Value = '  +4.82%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '5.16'
$c.Style = "Normal"
$ws.Range('E23').Value = '  +6.39%  '
$ws.Range('E24').Value = '  +0.17%  '
$ws.Range('E25').Value = '  +7.79%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '2.42'
$c.Style = "Normal"
$ws.Range('E26').Value = '  +0.19%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '166.45'
$c.Style = "Normal"
$ws.Range('E27').Value = '  -2.49%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '8.80'
$c.Style = "Normal"
$ws.Range('E28').Value = '  +4.15%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '18.78'
$c.Style = "Normal"
$ws.Range('E29').Value = '  +1.98%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '0.132'
$c.Style = "Normal"
$ws.Range('E30').Value = '  +4.51%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '4.128.88'
$c.Style = "Normal"
$ws.Range('E31').Value = '  -0.82%  '
$ws.Range('B32').Value = 'TrustWalletToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '1.68'
$c.Style = "Normal"
$ws.Range('E32').Value = '  +25.67%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '4.38'
$c.Style = "Normal"
$ws.Range('E33').Value = '  +5.46%  '
$ws.Range('E34').Value = '  +14.73%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '0.0592'
$c.Style = "Normal"
$ws.Range('E35').Value = '  +4.70%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '4.27'
$c.Style = "Normal"
$ws.Range('E36').Value = '  +4.22%  '
$ws.Range('E37').Value = '  +0.46%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '0.916'
$c.Style = "Normal"
$ws.Range('E38').Value = '  -3.39%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '2.05'
$c.Style = "Normal"
$ws.Range('E39').Value = '  -0.51%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '99.60'
$c.Style = "Normal"
$ws.Range('E40').Value = '  +10.64%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '0.0221'
$c.Style = "Normal"
$ws.Range('E41').Value = '  +5.77%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '17.03'
$c.Style = "Normal"
$ws.Range('E42').Value = '  +4.40%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '1.13'
$c.Style = "Normal"
$ws.Range('E43').Value = '  +1.41%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '0.0652'
$c.Style = "Normal"
$ws.Range('E44').Value = '  +0.74%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '1.344.75'
$c.Style = "Normal"
$ws.Range('E45').Value = '  +0.54%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '2.46'
$c.Style = "Normal"
$ws.Range('E46').Value = '  +2.10%  '
$ws.Range('E47').Value = '  +1.61%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '6.74'
$c.Style = "Normal"
$ws.Range('E48').Value = '  +3.64%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '2.76'
$c.Style = "Normal"
$ws.Range('E49').Value = '  -1.41%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '44.99'
$c.Style = "Normal"
$ws.Range('E50').Value = '  -9.05%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '0.0759'
$c.Style = "Normal"
$ws.Range('E51').Value = '  +7.10%  '
